$wb = $excel.ActiveWorkbook

# --- Sheet "Oct-24 RT Pk" (Peak) ---
$wsPk = $wb.Worksheets.Item("Oct-24 RT Pk")

# O3 note text changes
$wsPk.Range("O3").Value = "Summer ratings. Eastward wind transfer from MISO into SPP. Lisbon-Enderlin 69 kV outage."

# New notes added in row 4
$wsPk.Range("O4").Value = "Southward wind transfer from the greater Sioux City, IA area (including SPP) into Omaha, NE. Flows show a strong negative correlation with Walter Scott Jr. output."
$wsPk.Range("P4").Value = "Curious B-rating uprate on 11/8/2024 from 217 MW to 284 MW. This was followed by a 1/30/2025 B-rating derate from 284 MW to 259 MW then again on 2/21/2025 from 259 MW to 256 MW and again on 4/7/2025 from 256 MW to 217 MW. Tekamah-Oakland 115 kV is a switching solution that relieves Tekamah-Substation 1226 at the expense of increasing flows on Raun-Tekamah."

$wsPk.Select()
$excel.ActiveWindow.ScrollColumn = 2
$wsPk.Range("O12").Select()

# --- Sheet "Oct-24 RT Off" (Off-peak) ---
$wsOff = $wb.Worksheets.Item("Oct-24 RT Off")

# New "notes" header column (copy header formatting from O1, then set the text)
$wsOff.Range("O1").Copy()
$wsOff.Range("P1").PasteSpecial(-4122)
$wsOff.Range("P1").Value = "notes"

# New note content
$wsOff.Range("O2").Value = "Summer ratings. Eastward wind transfer from MISO into SPP. Lisbon-Enderlin 69 kV outage."
$wsOff.Range("P2").Value = "Winter ratings from 11/1-4/1. Winter (only) B-ratings uprated from 131 MW to 169 MW in Winter 2024-25."

$wsOff.Range("O4").Value = "Southward wind transfer from the greater Sioux City, IA area (including SPP) into Omaha, NE. Flows show a strong negative correlation with Walter Scott Jr. output."
$wsOff.Range("P4").Value = "Curious B-rating uprate on 11/8/2024 from 217 MW to 284 MW. This was followed by a 1/30/2025 B-rating derate from 284 MW to 259 MW then again on 2/21/2025 from 259 MW to 256 MW and again on 4/7/2025 from 256 MW to 217 MW. Tekamah-Oakland 115 kV is a switching solution that relieves Tekamah-Substation 1226 at the expense of increasing flows on Raun-Tekamah."

$wsOff.Select()
$excel.ActiveWindow.ScrollColumn = 2
$wsOff.Range("O4:P4").Select()

# Restore the originally active sheet/tab
$wsPk.Select()
